# Agrega la cuenta 40139, ENCINA FERRETERIA, como nueva fila 195
# (desplaza las filas 195:322 a 196:323).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja 1")
$ws.Activate()

# Inserta una fila nueva en la posicion 195, empujando el resto hacia abajo
$ws.Rows.Item(195).Insert()

# Completa los datos de la nueva fila: CUENTA, NOMBRE, CATEGORIAS, LISTA_PRECIOS
$ws.Cells.Item(195, 1).Value = 40139
$ws.Cells.Item(195, 2).Value = "ENCINA FERRETERIA"
$ws.Cells.Item(195, 3).Value = "SANITARIOS,REPUESTOS,BULONERIA,PINTURAS,FERRETERIA,ELECTRICIDAD"
$ws.Cells.Item(195, 4).Value = "E"

# Da formato a C195:D195 (borde medio gris, igual al resto de la tabla)
$rngFmt = $ws.Range("C195:D195")
$rngFmt.Borders.LineStyle = 1
$rngFmt.Borders.Weight = -4138
$rngFmt.Borders.Color = 13421772

# Restaura la altura de fila estandar usada en el resto de la hoja
$ws.Rows.Item(195).RowHeight = 15.75

# Actualiza el autofiltro para que cubra la nueva ultima fila (D323)
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A1:D323").AutoFilter()

# Actualiza el nombre definido oculto _FilterDatabase generado por el autofiltro
$wb.Names.Item("Hoja 1!_FilterDatabase").RefersTo = "='Hoja 1'!`$A`$1:`$D`$323"

# Restaura la seleccion activa tal como queda en el archivo final
$ws.Range("C195:D195").Select()
